# Auto-generated PowerShell COM-interop script to append new SeniorConnect
# sensor-log rows (2026-01-28, ~17:12-17:18) to the end of several sheets.
#
# Column layout on every sheet: A=Date, B=Timestamp, C=Hour, D=Location,
# E=Value, F=Status. Column A (dates) is forced to Text format before
# assignment so Excel does not silently reinterpret "2026-01-28" as a date
# serial number; likewise the Humidity sheet's E column (percentages like
# "87.6%") is forced to Text so it is not reinterpreted as a numeric ratio.

$wb = $excel.ActiveWorkbook

# --- PIR sheet: rows 125-139 (Bathroom "No Motion" / "Inactive") ---
$ws = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-01-28", "17:12:22", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:12:23", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:12:25", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:17:27", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:17:32", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:17:37", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:17:42", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:17:47", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:17:53", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:17:57", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:18:02", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:18:07", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:18:13", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:18:18", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:18:23", "17:00", "Bathroom", "No Motion", "Inactive")
)
$r = 125
foreach ($row in $pirRows) {
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $r = $r + 1
}

# --- Humidity sheet: rows 126-138 (Bathroom "%" / "Active") ---
$ws = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-01-28", "17:12:22", "17:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-28", "17:12:24", "17:00", "Bathroom", "86.7%", "Active"),
    @("2026-01-28", "17:12:26", "17:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-28", "17:17:28", "17:00", "Bathroom", "86.6%", "Active"),
    @("2026-01-28", "17:17:32", "17:00", "Bathroom", "87.5%", "Active"),
    @("2026-01-28", "17:17:36", "17:00", "Bathroom", "87.5%", "Active"),
    @("2026-01-28", "17:17:40", "17:00", "Bathroom", "86.6%", "Active"),
    @("2026-01-28", "17:17:48", "17:00", "Bathroom", "86.5%", "Active"),
    @("2026-01-28", "17:17:52", "17:00", "Bathroom", "86.2%", "Active"),
    @("2026-01-28", "17:18:00", "17:00", "Bathroom", "86.5%", "Active"),
    @("2026-01-28", "17:18:12", "17:00", "Bathroom", "86.5%", "Active"),
    @("2026-01-28", "17:18:16", "17:00", "Bathroom", "87.5%", "Active"),
    @("2026-01-28", "17:18:24", "17:00", "Bathroom", "87.4%", "Active")
)
$r = 126
foreach ($row in $humidityRows) {
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $r = $r + 1
}

# --- Temperature sheet: rows 125-138 (Bathroom "...C" / "Active") ---
$ws = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-01-28", "17:12:21", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:12:23", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:12:24", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:12:26", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:17:28", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:17:32", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:17:36", "17:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "17:17:40", "17:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "17:17:48", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:17:52", "17:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "17:18:00", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:18:12", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:18:16", "17:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "17:18:24", "17:00", "Bathroom", "22.8C", "Active")
)
$r = 125
foreach ($row in $temperatureRows) {
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $r = $r + 1
}

# --- Proximity sheet: row 3 (Living Room Main Door entry event) ---
$ws = $wb.Worksheets.Item("Proximity")
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2026-01-28"
$ws.Range("B3").Value = "17:18:19"
$ws.Range("C3").Value = "17:00"
$ws.Range("D3").Value = "Living Room Main Door"
$ws.Range("E3").Value = "ENTER"
$ws.Range("F3").Value = "User ENTERED Living Room Main Door"

# --- Camera sheet: row 3 (Living Room Main Door image captured) ---
$ws = $wb.Worksheets.Item("Camera")
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2026-01-28"
$ws.Range("B3").Value = "17:18:20"
$ws.Range("C3").Value = "17:00"
$ws.Range("D3").Value = "Living Room Main Door"
$ws.Range("E3").Value = "Image Captured"
$ws.Range("F3").Value = "Active"
